$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: 四方坪站
$ws.Cells.Item(16, 1).Value = 46030
$ws.Cells.Item(16, 2).Value = "四方坪站"
$ws.Cells.Item(16, 3).Value = 14284.84
$ws.Cells.Item(16, 4).Value = 10035.030000000001
$ws.Cells.Item(16, 5).Value = 3199.13
$ws.Cells.Item(16, 6).Value = 598

# New row 17: 高岭站
$ws.Cells.Item(17, 1).Value = 46030
$ws.Cells.Item(17, 2).Value = "高岭站"
$ws.Cells.Item(17, 3).Value = 4706.1899999999996
$ws.Cells.Item(17, 4).Value = 4036.5
$ws.Cells.Item(17, 5).Value = 1208.8
$ws.Cells.Item(17, 6).Value = 174

$ws.Cells.Item(18, 9).Select()
